$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) Merge the split runs back into single runs (no visible text change) ---

$shp1 = $s.Shapes.Item(8).GroupItems.Item(2)
$shp1.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Swagger Reference Properties"

$shp2 = $s.Shapes.Item(14)
$shp2.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "Cycling: Repeat key combination to cycle file, project & workspace scope"

# --- 2) Add the new straight-line annotation connector ---

# A throw-away shape "burns" id=2 so the real connector lands on id=3,
# matching the id PowerPoint assigned in the authored deck.
$burn = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$burn.Delete()

$conn = $s.Shapes.AddLine(432.0, 373.0793, 432.0, 498.0)
$conn.Name = "Straight Connector 2"

$conn.Line.ForeColor.RGB = 10498160
$conn.Shadow.Type = 21
$conn.Shadow.Visible = $true
$conn.Shadow.Style = 2
$conn.Shadow.Blur = 4
$conn.Shadow.OffsetX = 0
$conn.Shadow.OffsetY = 3
$conn.Shadow.ForeColor.RGB = 0
$conn.Shadow.Transparency = 0.6

Write-Output $conn.Name
